# repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to match the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 0
$ws.Range("F5").Value  = -5
$ws.Range("F7").Value  = -9
$ws.Range("F8").Value  = -5
$ws.Range("F11").Value = 9
$ws.Range("F13").Value = 4
$ws.Range("F20").Value = 5
